$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-09-04 14:48:57"

$wsZhCn.Range("H3").Value = "2016-09-04 14:48:52"
$wsZhCn.Range("K3").Value = "2016-09-04 14:49:33"

$wsDeDe.Range("K3").Value = "2016-09-04 14:49:40"
